# clinEpi_eda_inconsistencyIssues.xlsx -- "some inconsistency issue tracker updates"
#
# Updates to sheet "termWithDifferentLabels" (comments column E), sheet
# "LabelsUsedMultipleTerms" (comment column F), sheet "termWithDifferentParent"
# (parent-category column F), and sheet "units_issues" (comment column E).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: termWithDifferentLabels
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("termWithDifferentLabels")

# NOTE: assignment order below matches the order brand-new comment strings
# were introduced (so freshly minted shared-string entries land at the same
# indices the saved workbook used). Cells that reuse an already-existing
# string (E90, E96, E101) are interleaved wherever convenient.
$ws1.Range("E93").Value  = "Fever at episode >37.5 C"
$ws1.Range("E94").Value  = "Persistent diarrheal episode (>14 days)"
$ws1.Range("E91").Value  = "Maximum loose stools at episode count"
$ws1.Range("E97").Value  = "Severe anemia (hemoglobin <5 mg/dL)"
$ws1.Range("E98").Value  = "Persons enrolled in study count"
$ws1.Range("E106").Value = "Persons <=18 years living in house"
$ws1.Range("E102").Value = "Eukaryota in stool"
$ws1.Range("E104").Value = "Eukaryota in urine"
$ws1.Range("E103").Value = "Schistosoma haematobium infection intensity, by microscopy"

$ws1.Range("E90").Value  = "Children <60 months in dwelling count"
$ws1.Range("E96").Value  = "Diarrhea treatment"
$ws1.Range("E101").Value = "Schistosoma mansoni infection intensity, by microscopy"

# view state: active cell moved, tab selection moved off this sheet
$ws1.Range("A1").Select()

# ---------------------------------------------------------------------------
# Sheet 2: LabelsUsedMultipleTerms
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LabelsUsedMultipleTerms")

$ws2.Range("F33").Value = "Use EUPATH_0033152, should be under Breastfeeding. Remove Breastfeeding summary category"

# ---------------------------------------------------------------------------
# Sheet 3: termWithDifferentParent
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("termWithDifferentParent")

$ws3.Range("F120").Value = "Shouldn't exist anymore. All variables should be moved under Breastfeeding EUPATH_0011730 instead"
$ws3.Range("F136").Value = "Village observation details"

$ws3.Range("F118").Value = "Administrative information"
$ws3.Range("F123").Value = "Symptoms"
$ws3.Range("F124").Value = "Symptoms"
$ws3.Range("F127").Value = "Symptoms"
$ws3.Range("F130").Value = "Symptoms"
$ws3.Range("F131").Value = "Symptoms"
$ws3.Range("F133").Value = "Symptoms"
$ws3.Range("F135").Value = "Symptoms"
$ws3.Range("F137").Value = "Symptoms"

# ---------------------------------------------------------------------------
# Sheet 4: units_issues
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("units_issues")

$ws4.Range("E36").Value = "remove unit"
# Column E has no sheet-level default style in this sheet, so the new cell
# needs its format copied explicitly from a neighboring comment cell (E1)
# that already carries the correct style (matches the other "comments" cells).
$ws4.Range("E1").Copy()
$ws4.Range("E36").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Window / view state changes
# ---------------------------------------------------------------------------
# Active tab moved from sheet 1 (termWithDifferentLabels) to sheet 2
# (LabelsUsedMultipleTerms).
$ws2.Activate()
$ws2.Range("A39").Select()

$excel.ActiveWindow.ScrollRow = 26
